# Update res_bus/vm_pu.xlsx results for the "case with 380 kV" run.
# Rows 2-25 hold the per-bus voltage magnitudes (p.u.); column G (the
# reference bus) is unchanged and stays at 1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$columns = @("B","C","D","E","F","I","J","K","L","M","N")

# Each row: row number followed by the new values for $columns, in order.
$newValues = @(
    @(2, 1.02, 1.015095992252911, 1.040238998888775, 1.016781028109505, 1.043223671430217, 1.033401191578782, 1.020323464360613, 1.043021912275336, 1.019632176276497, 1.045998141876539, 1.021772440712743),
    @(3, 1.02, 1.016195744974358, 1.040819828612072, 1.017718670370386, 1.044070901047487, 1.033479741128522, 1.021057288576413, 1.043413729032948, 1.020374906376732, 1.046656268083197, 1.022507307043099),
    @(4, 1.02, 1.016907851338674, 1.041195370820844, 1.018326160095925, 1.044619284223505, 1.033529078007704, 1.021532111576009, 1.04366626757097, 1.020855657432625, 1.04708164529066, 1.022982804345818),
    @(5, 1.02, 1.017207339219309, 1.041353176672695, 1.018581733313434, 1.044849863291978, 1.033549461862354, 1.02173172451236, 1.043772195655157, 1.0210578020243, 1.047260358832596, 1.023182700755435),
    @(6, 1.02, 1.017257631504594, 1.041379668688471, 1.018624655991075, 1.044888580757741, 1.033552863425224, 1.021765240242057, 1.043789967384828, 1.021091745143897, 1.047290358834647, 1.023216264081311),
    @(7, 1.02, 1.016911852644805, 1.041197479715969, 1.018329574353158, 1.044622365082036, 1.033529351782838, 1.02153477882448, 1.043667683928518, 1.020858358354096, 1.04708403372361, 1.022985475382088),
    @(8, 1.02, 1.015467557022683, 1.040435352673367, 1.017097748513982, 1.043509960659313, 1.033428045860716, 1.020571465722849, 1.04315453339171, 1.019883153135218, 1.046220656077353, 1.022020794265359),
    @(9, 1.02, 1.012926280662738, 1.039090220769406, 1.014933047619323, 1.041551129408298, 1.033238150575161, 1.018873913405392, 1.042242750693426, 1.018165915589622, 1.044695704391946, 1.020320831228899),
    @(10, 1.02, 1.011234595853459, 1.038192118091489, 1.013493924453952, 1.04024625936632, 1.033103941531001, 1.017742172099673, 1.041629902726324, 1.017021910037742, 1.043676754844614, 1.019187482720713),
    @(11, 1.02, 1.010502663799955, 1.037802932783046, 1.012871723379454, 1.03968149793555, 1.033044030370619, 1.017252107312908, 1.041363365409706, 1.016526739639982, 1.04323500700686, 1.018696721985741),
    @(12, 1.02, 1.010230877909102, 1.037658328647601, 1.012640752971669, 1.039471760449758, 1.033021507252932, 1.017070073438113, 1.041264187121252, 1.016342840195048, 1.043070843055642, 1.018514429601963),
    @(13, 1.02, 1.010289173012392, 1.037689348659469, 1.012690290443526, 1.039516748014055, 1.033026350720713, 1.017109120412174, 1.041285469074194, 1.016382285941217, 1.043106060339275, 1.018553532027206),
    @(14, 1.02, 1.010480196137541, 1.037790980634576, 1.012852628379562, 1.039664160143141, 1.033042174095546, 1.017237060376616, 1.041355170852716, 1.016511537864544, 1.043221438769779, 1.018681653681073),
    @(15, 1.02, 1.010597903200283, 1.037853593754488, 1.012952669065908, 1.039754990976865, 1.03305188771426, 1.017315888147544, 1.04139809334085, 1.01659117807164, 1.043292516783957, 1.018760593396477),
    @(16, 1.02, 1.011283183907078, 1.038217940854692, 1.013535237899274, 1.040283746251562, 1.033107879816026, 1.017774695805485, 1.041647567391419, 1.01705477691136, 1.043706061020134, 1.019220052613918),
    @(17, 1.02, 1.011713196831239, 1.038446406891513, 1.013900922110954, 1.040615490116658, 1.033142521328997, 1.018062490164747, 1.041803743545598, 1.0173456315768, 1.043965323882387, 1.019508255674181),
    @(18, 1.02, 1.011964071968862, 1.038579638108731, 1.014114311049735, 1.040809015284306, 1.03316255363644, 1.018230354458122, 1.04189472546697, 1.017515300657155, 1.044116495718447, 1.019676358354107),
    @(19, 1.02, 1.012049623405305, 1.038625061509622, 1.014187086680493, 1.040875006486187, 1.033169354687089, 1.018287591625544, 1.041925728745772, 1.017573156559286, 1.044168032578633, 1.019733676804872),
    @(20, 1.02, 1.011667054739444, 1.038421897650539, 1.013861678187195, 1.040579894585041, 1.033138822563831, 1.018031612686646, 1.041786998999416, 1.017314423712007, 1.04393751277929, 1.01947733434652),
    @(21, 1.02, 1.010423942212374, 1.037761053733981, 1.012804819960383, 1.039620749826652, 1.033037521938858, 1.017199385300824, 1.041334650221524, 1.016473475596299, 1.04318746487114, 1.018643925102352),
    @(22, 1.02, 1.009642846792458, 1.037345304829998, 1.01214115741557, 1.039017929990805, 1.032972271607938, 1.016676119389346, 1.041049231866688, 1.015944905333075, 1.042715422811725, 1.01811991609326),
    @(23, 1.02, 1.010056873198281, 1.037565724289506, 1.012492899097114, 1.039337473613448, 1.033007009565112, 1.016953513637714, 1.041200632708639, 1.016225094545404, 1.042965704111471, 1.018397704273278),
    @(24, 1.02, 1.011687904191659, 1.038432972411834, 1.013879410544231, 1.04059597859776, 1.03314049441275, 1.018045564893762, 1.041794565488315, 1.017328525145606, 1.04395007957984, 1.019491306367369),
    @(25, 1.02, 1.013582819047461, 1.039438215143753, 1.015491969331991, 1.042057362067115, 1.033288587677214, 1.019312778996914, 1.04247935427774, 1.018609719219104, 1.045090355223863, 1.020760320059899),
)

foreach ($entry in $newValues) {
    $rowNum = $entry[0]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $ws.Range("$($columns[$i])$rowNum").Value = $entry[$i + 1]
    }
}